$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 01:22"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 309728
$ws.Range("C4").Value = 32567
$ws.Range("E4").Value = 286546
$ws.Range("G4").Value = 1037
$ws.Range("H4").Value = 8441

# --- España (row 5) ---
$ws.Range("F5").Value = 6532

# --- Brasil (row 19) ---
$ws.Range("B19").Value = 10360
$ws.Range("C19").Value = 1166
$ws.Range("E19").Value = 9788
$ws.Range("G19").Value = 82
$ws.Range("H19").Value = 445

# --- Rows 34-36: Japon / Filipinas / India re-sorted by updated totals ---
# Row 34 becomes Japon with freshly updated figures
$ws.Range("A34").Value = "Japon"
$ws.Range("B34").Value = 3139
$ws.Range("C34").Value = 204
$ws.Range("D34").Value = 514
$ws.Range("E34").Value = 2548
$ws.Range("F34").Value = 64
$ws.Range("G34").Value = 8
$ws.Range("H34").Value = 77

# Row 35 becomes Filipinas (previous Filipinas figures)
$ws.Range("A35").Value = "Filipinas"
$ws.Range("B35").Value = 3094
$ws.Range("C35").Value = 76
$ws.Range("D35").Value = 57
$ws.Range("E35").Value = 2893
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 8
$ws.Range("H35").Value = 144

# Row 36 becomes India (previous India figures)
$ws.Range("A36").Value = "India"
$ws.Range("B36").Value = 3082
$ws.Range("C36").Value = 23
$ws.Range("D36").Value = 229
$ws.Range("E36").Value = 2767
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 86

# --- San Marino (row 98) ---
$ws.Range("B98").Value = 259
$ws.Range("C98").Value = 8
$ws.Range("D98").Value = 27
$ws.Range("E98").Value = 200
$ws.Range("F98").Value = 14

# --- Rows 133-134: Guatemala / Guayana Francesa re-sorted by updated totals ---
# Row 133 becomes Guatemala with freshly updated figures
$ws.Range("A133").Value = "Guatemala"
$ws.Range("B133").Value = 61
$ws.Range("C133").Value = 11
$ws.Range("D133").Value = 15
$ws.Range("E133").Value = 44
$ws.Range("F133").Value = 1
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 2

# Row 134 becomes Guayana Francesa (previous Guayana Francesa figures)
$ws.Range("A134").Value = "Guayana Francesa"
$ws.Range("B134").Value = 61
$ws.Range("C134").Value = 4
$ws.Range("D134").Value = 22
$ws.Range("E134").Value = 39
$ws.Range("F134").Value = 1
$ws.Range("G134").Value = 0
$ws.Range("H134").Value = 0
